# --- Workbook-level changes -------------------------------------------------
$wb = $excel.ActiveWorkbook

# Rename "ex3g" -> "ex3"
$ex3 = $wb.Worksheets.Item("ex3g")
$ex3.Name = "ex3"

# Add a new trailing sheet named "trash" (placed after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$trash = $wb.Worksheets.Add($null, $lastSheet)
$trash.Name = "trash"

# --- ex2g: selection only change --------------------------------------------
$ws2 = $wb.Worksheets.Item("ex2g")
$ws2.Activate()
$ws2.Range("F16").Select()

# --- ex3: values/formula changes + new IRR row -------------------------------
$ex3.Activate()
$ex3.Range("B10").Value = -440000
$ex3.Range("F10").Formula = "=SUM(B10:D10)"

# --- ex5g: selection only change ---------------------------------------------
$ws5 = $wb.Worksheets.Item("ex5g")
$ws5.Activate()
$ws5.Range("F10").Select()

# --- ex9g: selection only change (no longer the active tab) ------------------
$ws9 = $wb.Worksheets.Item("ex9g")
$ws9.Activate()
$ws9.Range("G10").Select()

# --- trash: new sheet contents ------------------------------------------------
$trash.Activate()
$trash.Application.ActiveWindow.Zoom = 220

$trash.Range("B2").Value = 0.02
$trash.Range("B2").NumberFormat = "0%"

$trash.Range("B3").Formula = "=1*(1+B2)^(12)"
$trash.Range("C3").Formula = "=(B3*100)/1"

$trash.Range("C4").Formula = "=(C3-100)/100"
$trash.Range("D4").Formula = "=C4*12"

$trash.Range("D5").Formula = "=1/12"

$trash.Range("B6").Formula = "=1*(1+C4)^(0.08333)"

$trash.Range("B8").Formula = "=1*(1+0.168)^(1/360)"
$trash.Range("D8").Value = 32050
$trash.Range("E8").Formula = "=(D8*100)/40"

# NOTE: "asumo" must be written before "TIR" so the shared-string table
# picks up the same ordering as the target workbook.
$trash.Range("A9").Value = "asumo"
$trash.Range("B9").Value = 0.00043162597291044902
$trash.Range("B9").NumberFormat = "0%"
$trash.Range("E9").Formula = "=E8/(1+0.000431)^(108)"
$trash.Range("F9").Formula = "=(25*E8)/100"

$trash.Range("B10").Formula = "=1*(1+B9)^(360)"
$trash.Range("F10").Formula = "=F9*(1+0.000431)^(49)"

$trash.Range("B12").Formula = "=(35*E9)/100"

$trash.Range("B13").Formula = "=B12+F10+D8"

$trash.Range("B14").Formula = "=E9-B13"

$trash.Range("B15").Select()

# --- ex3: add the new "TIR" row now, after "asumo" already exists in sst -----
$ex3.Activate()
$ex3.Range("C13").Value = "TIR"
$ex3.Range("D13").Formula = "=IRR(B10:D10)"
$ex3.Range("D13").NumberFormat = "0%"
$ex3.Range("H10").Select()

# Final state: "ex3" is the active tab (index 2)
$ex3.Activate()
